$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 59
$ws.Range("B19").Value = 'How many tracks can you define in one ODF?'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = '20.'
$ws.Range("C19").ClearFormats() | Out-Null
$ws.Range("D19").Value = $true
$ws.Range("E19").Value = '200 tracks can be defined in one ODF.'
$ws.Range("F19").Value = 0.230478972196579

$ws.Range("A20").Value = 60
$ws.Range("B20").Value = 'How many tracks can you define in one ODF?'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = '20.'
$ws.Range("C20").ClearFormats() | Out-Null
$ws.Range("D20").Value = $true
$ws.Range("E20").Value = '200 tracks can be defined in one ODF.'
$ws.Range("F20").Value = 0.230478972196579

$ws.Range("A21").Value = 61
$ws.Range("B21").Value = 'How many curve shades can I create?'
$ws.Range("C21").Value = 'According to the document, you can create 250 curve shades.'
$ws.Range("D21").Value = $true
$ws.Range("E21").Value = '250 curve shades can be created.'
$ws.Range("F21").Value = 0.9534614086151123

$ws.Range("A22").Value = 62
$ws.Range("B22").Value = 'How many curve shades can I create?'
$ws.Range("C22").Value = 'According to the document, you can create 250 curve shades.'
$ws.Range("D22").Value = $true
$ws.Range("E22").Value = '250 curve shades can be created.'
$ws.Range("F22").Value = 0.9534614086151123

$ws.Range("A23").Value = 63
$ws.Range("B23").Value = 'How many curves can I load in one go?'
$ws.Range("C23").Value = 'According to the GEO application documentation, you can load up to 5 data files to form one curve.'
$ws.Range("D23").Value = $true
$ws.Range("E23").Value = '450 curves can be loaded in one go.'
$ws.Range("F23").Value = 0.647385835647583

$ws.Range("A24").Value = 64
$ws.Range("B24").Value = 'How many curves can I load in one go?'
$ws.Range("C24").Value = 'According to the GEO application documentation, you can load up to 5 data files to form one curve.'
$ws.Range("D24").Value = $true
$ws.Range("E24").Value = '450 curves can be loaded in one go.'
$ws.Range("F24").Value = 0.647385835647583

$ws.Range("A25").Value = 67
$ws.Range("B25").Value = 'How many tables can I have in my log?'
$ws.Range("C25").Value = 'You can have up to 100 tables in a log.'
$ws.Range("D25").Value = $true
$ws.Range("E25").Value = '100 tables can be presented in a log.'
$ws.Range("F25").Value = 0.8799165487289429

$ws.Range("A26").Value = 68
$ws.Range("B26").Value = 'How many tables can I have in my log?'
$ws.Range("C26").Value = 'You can have up to 100 tables in a log.'
$ws.Range("D26").Value = $true
$ws.Range("E26").Value = '100 tables can be presented in a log.'
$ws.Range("F26").Value = 0.8799165487289429

$ws.Range("A27").Value = 71
$ws.Range("B27").Value = 'How many symbols can I have in the plot at any one time?'
$ws.Range("C27").Value = 'According to the document, you can have 10000 symbols per plot.'
$ws.Range("D27").Value = $true
$ws.Range("E27").Value = '10000 symbols can be defined in a plot at any one time.'
$ws.Range("F27").Value = 0.8716347217559814

$ws.Range("A28").Value = 72
$ws.Range("B28").Value = 'How many symbols can I have in the plot at any one time?'
$ws.Range("C28").Value = 'According to the document, you can have 10000 symbols per plot.'
$ws.Range("D28").Value = $true
$ws.Range("E28").Value = '10000 symbols can be defined in a plot at any one time.'
$ws.Range("F28").Value = 0.8716347217559814

$ws.Range("A29").Value = 73
$ws.Range("B29").Value = 'How many scales can I define?'
$ws.Range("C29").Value = 'According to the document, you can define up to 23 scales.'
$ws.Range("D29").Value = $true
$ws.Range("E29").Value = '23 scales can be defined.'
$ws.Range("F29").Value = 0.8617551922798157

$ws.Range("A30").Value = 74
$ws.Range("B30").Value = 'How many scales can I define?'
$ws.Range("C30").Value = 'According to the document, you can define up to 23 scales.'
$ws.Range("D30").Value = $true
$ws.Range("E30").Value = '23 scales can be defined.'
$ws.Range("F30").Value = 0.8617551922798157

$ws.Range("A31").Value = 79
$ws.Range("B31").Value = 'How many tracks can you define in one ODF?'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = '20.'
$ws.Range("C31").ClearFormats() | Out-Null
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = '200 tracks can be defined in one ODF.'
$ws.Range("F31").Value = 0.230478972196579

$ws.Range("A32").Value = 80
$ws.Range("B32").Value = 'How many tracks can you define in one ODF?'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = '20.'
$ws.Range("C32").ClearFormats() | Out-Null
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = '200 tracks can be defined in one ODF.'
$ws.Range("F32").Value = 0.230478972196579

$ws.Range("A33").Value = 81
$ws.Range("B33").Value = 'How many curve shades can I create?'
$ws.Range("C33").Value = 'According to the document, you can create 250 curve shades.'
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = '250 curve shades can be created.'
$ws.Range("F33").Value = 0.9534614086151123

$ws.Range("A34").Value = 82
$ws.Range("B34").Value = 'How many curve shades can I create?'
$ws.Range("C34").Value = 'According to the document, you can create 250 curve shades.'
$ws.Range("D34").Value = $true
$ws.Range("E34").Value = '250 curve shades can be created.'
$ws.Range("F34").Value = 0.9534614086151123

$ws.Range("A35").Value = 83
$ws.Range("B35").Value = 'How many curves can I load in one go?'
$ws.Range("C35").Value = 'According to the GEO application documentation, you can load up to 5 data files to form one curve.'
$ws.Range("D35").Value = $true
$ws.Range("E35").Value = '450 curves can be loaded in one go.'
$ws.Range("F35").Value = 0.647385835647583

$ws.Range("A36").Value = 84
$ws.Range("B36").Value = 'How many curves can I load in one go?'
$ws.Range("C36").Value = 'According to the GEO application documentation, you can load up to 5 data files to form one curve.'
$ws.Range("D36").Value = $true
$ws.Range("E36").Value = '450 curves can be loaded in one go.'
$ws.Range("F36").Value = 0.647385835647583

$ws.Range("A37").Value = 87
$ws.Range("B37").Value = 'How many tables can I have in my log?'
$ws.Range("C37").Value = 'You can have up to 100 tables in a log.'
$ws.Range("D37").Value = $true
$ws.Range("E37").Value = '100 tables can be presented in a log.'
$ws.Range("F37").Value = 0.8799165487289429

$ws.Range("A38").Value = 88
$ws.Range("B38").Value = 'How many tables can I have in my log?'
$ws.Range("C38").Value = 'You can have up to 100 tables in a log.'
$ws.Range("D38").Value = $true
$ws.Range("E38").Value = '100 tables can be presented in a log.'
$ws.Range("F38").Value = 0.8799165487289429

$ws.Range("A39").Value = 91
$ws.Range("B39").Value = 'How many symbols can I have in the plot at any one time?'
$ws.Range("C39").Value = 'According to the document, you can have 10000 symbols per plot.'
$ws.Range("D39").Value = $true
$ws.Range("E39").Value = '10000 symbols can be defined in a plot at any one time.'
$ws.Range("F39").Value = 0.8716347217559814

$ws.Range("A40").Value = 92
$ws.Range("B40").Value = 'How many symbols can I have in the plot at any one time?'
$ws.Range("C40").Value = 'According to the document, you can have 10000 symbols per plot.'
$ws.Range("D40").Value = $true
$ws.Range("E40").Value = '10000 symbols can be defined in a plot at any one time.'
$ws.Range("F40").Value = 0.8716347217559814

$ws.Range("A41").Value = 93
$ws.Range("B41").Value = 'How many scales can I define?'
$ws.Range("C41").Value = 'According to the document, you can define 23 scales.'
$ws.Range("D41").Value = $true
$ws.Range("E41").Value = '23 scales can be defined.'
$ws.Range("F41").Value = 0.9348680973052979

$ws.Range("A42").Value = 94
$ws.Range("B42").Value = 'How many scales can I define?'
$ws.Range("C42").Value = 'According to the document, you can define 23 scales.'
$ws.Range("D42").Value = $true
$ws.Range("E42").Value = '23 scales can be defined.'
$ws.Range("F42").Value = 0.9348680973052979

# Match the existing "A" column style (bold, bordered, centered/top-aligned)
# used throughout the table for the newly added rows.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19:A42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
